# Populate the contact_submissions worksheet with header row + two data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - existing cells already carry style s="1" (bold + border + centered)
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Department"
$ws.Range("D1").Value = "Year & Section"
$ws.Range("E1").Value = "Email"
$ws.Range("F1").Value = "Message"

# Row 2 data
$ws.Range("A2").Value = "2025-02-19 20:45:19"
$ws.Range("B2").Value = "Kawin"
$ws.Range("C2").Value = "k"
$ws.Range("D2").Value = "fgg"
$ws.Range("E2").Value = "kawin0275@gmail.com"
$ws.Range("F2").Value = "hi"

# Row 3 data
$ws.Range("A3").Value = "2025-02-19 20:45:44"
$ws.Range("B3").Value = "Kawin"
$ws.Range("C3").Value = "k"
$ws.Range("D3").Value = "fgg"
$ws.Range("E3").Value = "kawin0275@gmail.com"
$ws.Range("F3").Value = "ff"

# Selection moves to A1 (matches the target sheetView selection)
$ws.Range("A1").Select()
